$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2246.75
$ws.Cells.Item(19, 10).Value = 1995
$ws.Cells.Item(19, 12).Value = 1995
$ws.Cells.Item(19, 14).Value = -2345
$ws.Cells.Item(40, 8).Value = 103863.234
$ws.Cells.Item(40, 10).Value = 3999.875
$ws.Cells.Item(40, 12).Value = 3999.875
$ws.Cells.Item(40, 14).Value = -4349.875
$ws.Cells.Item(46, 8).Value = 8124.25
$ws.Cells.Item(46, 9).Value = 8332.333000000001
$ws.Cells.Item(46, 10).Value = 7500
$ws.Cells.Item(46, 11).Value = 24996.999
$ws.Cells.Item(46, 12).Value = 22500
$ws.Cells.Item(46, 13).Value = -24877.999
$ws.Cells.Item(46, 14).Value = -22738
$ws.Cells.Item(60, 8).Value = 8124.25
$ws.Cells.Item(60, 9).Value = 8332.333000000001
$ws.Cells.Item(60, 10).Value = 7500
$ws.Cells.Item(60, 11).Value = 24996.999
$ws.Cells.Item(60, 12).Value = 22500
$ws.Cells.Item(60, 13).Value = -24512.999
$ws.Cells.Item(60, 14).Value = -23468
$ws.Cells.Item(86, 8).Value = 132355960
$ws.Cells.Item(86, 9).Value = 76926504
$ws.Cells.Item(86, 10).Value = 312501700
$ws.Cells.Item(86, 11).Value = 76926504
$ws.Cells.Item(86, 12).Value = 312501700
$ws.Cells.Item(86, 13).Value = -76925381
$ws.Cells.Item(86, 14).Value = -312503946
$ws.Cells.Item(89, 8).Value = 132355960
$ws.Cells.Item(89, 9).Value = 76926504
$ws.Cells.Item(89, 10).Value = 312501700
$ws.Cells.Item(89, 11).Value = 384632520
$ws.Cells.Item(89, 12).Value = 1562508500
$ws.Cells.Item(89, 13).Value = -384626904
$ws.Cells.Item(89, 14).Value = -1562519732
$ws.Cells.Item(112, 8).Value = 4307308
$ws.Cells.Item(112, 10).Value = 4737329
$ws.Cells.Item(112, 12).Value = 14211987
$ws.Cells.Item(112, 14).Value = -14214203
$ws.Cells.Item(137, 8).Value = 21924.977
$ws.Cells.Item(137, 9).Value = 36990.793
$ws.Cells.Item(137, 11).Value = 110972.379
$ws.Cells.Item(137, 13).Value = -108422.379
$ws.Cells.Item(138, 8).Value = 2513.16
$ws.Cells.Item(138, 9).Value = 1150.9048
$ws.Cells.Item(138, 10).Value = 3499.6206
$ws.Cells.Item(138, 11).Value = 3452.7144
$ws.Cells.Item(138, 12).Value = 10498.8618
$ws.Cells.Item(138, 13).Value = 1687.2856
$ws.Cells.Item(138, 14).Value = -20778.8618

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18468734
$ws.Cells.Item(32, 9).Value = 17363342
$ws.Cells.Item(32, 11).Value = 17363342
$ws.Cells.Item(32, 13).Value = -17363055
$ws.Cells.Item(45, 8).Value = 5055
$ws.Cells.Item(45, 10).Value = 4831.6665
$ws.Cells.Item(45, 12).Value = 4831.6665
$ws.Cells.Item(45, 14).Value = -5585.6665
$ws.Cells.Item(61, 8).Value = 3079.077
$ws.Cells.Item(61, 9).Value = 2745.3044
$ws.Cells.Item(61, 11).Value = 2745.3044
$ws.Cells.Item(61, 13).Value = -2533.3044
$ws.Cells.Item(102, 8).Value = 1832.4
$ws.Cells.Item(102, 9).Value = 1646.125
$ws.Cells.Item(102, 10).Value = 2577.5
$ws.Cells.Item(102, 11).Value = 1646.125
$ws.Cells.Item(102, 12).Value = 2577.5
$ws.Cells.Item(102, 13).Value = -24.125
$ws.Cells.Item(102, 14).Value = -5821.5
$ws.Cells.Item(136, 8).Value = 3079.077
$ws.Cells.Item(136, 9).Value = 2745.3044
$ws.Cells.Item(136, 11).Value = 8235.913199999999
$ws.Cells.Item(136, 13).Value = -5685.913199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1456.6136
$ws.Cells.Item(86, 9).Value = 1311.2307
$ws.Cells.Item(86, 11).Value = 1311.2307
$ws.Cells.Item(86, 13).Value = -188.2307000000001
$ws.Cells.Item(89, 8).Value = 1456.6136
$ws.Cells.Item(89, 9).Value = 1311.2307
$ws.Cells.Item(89, 11).Value = 6556.1535
$ws.Cells.Item(89, 13).Value = -940.1535000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4539.2383
$ws.Cells.Item(31, 9).Value = 2972.0715
$ws.Cells.Item(31, 10).Value = 5322.8213
$ws.Cells.Item(31, 11).Value = 2972.0715
$ws.Cells.Item(31, 12).Value = 5322.8213
$ws.Cells.Item(31, 13).Value = -2677.0715
$ws.Cells.Item(31, 14).Value = -5912.8213
$ws.Cells.Item(34, 8).Value = 4539.2383
$ws.Cells.Item(34, 9).Value = 2972.0715
$ws.Cells.Item(34, 10).Value = 5322.8213
$ws.Cells.Item(34, 11).Value = 2972.0715
$ws.Cells.Item(34, 12).Value = 5322.8213
$ws.Cells.Item(34, 13).Value = -2770.0715
$ws.Cells.Item(34, 14).Value = -5726.8213
$ws.Cells.Item(58, 8).Value = 3566.6667
$ws.Cells.Item(58, 9).Value = 3197.28
$ws.Cells.Item(58, 11).Value = 3197.28
$ws.Cells.Item(58, 13).Value = -2994.28
$ws.Cells.Item(132, 8).Value = 3188.2375
$ws.Cells.Item(132, 9).Value = 3015.365
$ws.Cells.Item(132, 11).Value = 9046.094999999999
$ws.Cells.Item(132, 13).Value = -6516.094999999999
$ws.Cells.Item(134, 8).Value = 2580.25
$ws.Cells.Item(134, 9).Value = 2371.2122
$ws.Cells.Item(134, 11).Value = 7113.6366
$ws.Cells.Item(134, 13).Value = -4578.6366
$ws.Cells.Item(136, 8).Value = 3566.6667
$ws.Cells.Item(136, 9).Value = 3197.28
$ws.Cells.Item(136, 11).Value = 9591.84
$ws.Cells.Item(136, 13).Value = -7041.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 6778.778
$ws.Cells.Item(14, 9).Value = 6778.778
$ws.Cells.Item(14, 11).Value = 20336.334
$ws.Cells.Item(14, 13).Value = -20163.334
$ws.Cells.Item(46, 8).Value = 3163
$ws.Cells.Item(46, 10).Value = 5000
$ws.Cells.Item(46, 12).Value = 15000
$ws.Cells.Item(46, 14).Value = -15182
$ws.Cells.Item(86, 8).Value = 659.4
$ws.Cells.Item(86, 9).Value = 599
$ws.Cells.Item(86, 10).Value = 750
$ws.Cells.Item(86, 11).Value = 1797
$ws.Cells.Item(86, 12).Value = 2250
$ws.Cells.Item(86, 13).Value = -611
$ws.Cells.Item(86, 14).Value = -4622
$ws.Cells.Item(89, 8).Value = 659.4
$ws.Cells.Item(89, 9).Value = 599
$ws.Cells.Item(89, 10).Value = 750
$ws.Cells.Item(89, 11).Value = 5391
$ws.Cells.Item(89, 12).Value = 6750
$ws.Cells.Item(89, 13).Value = 537
$ws.Cells.Item(89, 14).Value = -18606
$ws.Cells.Item(122, 8).Value = 2629.4707
$ws.Cells.Item(122, 9).Value = 1599.5714
$ws.Cells.Item(122, 10).Value = 3350.4
$ws.Cells.Item(122, 11).Value = 14396.1426
$ws.Cells.Item(122, 12).Value = 30153.6
$ws.Cells.Item(122, 13).Value = -11946.1426
$ws.Cells.Item(122, 14).Value = -35053.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4488.2144
$ws.Cells.Item(70, 9).Value = 4493.25
$ws.Cells.Item(70, 10).Value = 4487.375
$ws.Cells.Item(70, 11).Value = 4493.25
$ws.Cells.Item(70, 12).Value = 4487.375
$ws.Cells.Item(70, 13).Value = -4223.25
$ws.Cells.Item(70, 14).Value = -5027.375
$ws.Cells.Item(73, 8).Value = 4488.2144
$ws.Cells.Item(73, 9).Value = 4493.25
$ws.Cells.Item(73, 10).Value = 4487.375
$ws.Cells.Item(73, 11).Value = 4493.25
$ws.Cells.Item(73, 12).Value = 4487.375
$ws.Cells.Item(73, 13).Value = -3557.25
$ws.Cells.Item(73, 14).Value = -6359.375
$ws.Cells.Item(100, 8).Value = 49799.8
$ws.Cells.Item(100, 10).Value = 49799.8
$ws.Cells.Item(100, 12).Value = 49799.8
$ws.Cells.Item(100, 14).Value = -51963.8
$ws.Cells.Item(102, 8).Value = 2858.7144
$ws.Cells.Item(102, 9).Value = 2875.5
$ws.Cells.Item(102, 11).Value = 2875.5
$ws.Cells.Item(102, 13).Value = -1253.5
$ws.Cells.Item(126, 8).Value = 3232.6667
$ws.Cells.Item(126, 9).Value = 3499
$ws.Cells.Item(126, 10).Value = 3099.5
$ws.Cells.Item(126, 11).Value = 10497
$ws.Cells.Item(126, 12).Value = 9298.5
$ws.Cells.Item(126, 13).Value = -8027
$ws.Cells.Item(126, 14).Value = -14238.5
$ws.Cells.Item(132, 8).Value = 4854.5806
$ws.Cells.Item(132, 9).Value = 5065.857
$ws.Cells.Item(132, 10).Value = 4410.9
$ws.Cells.Item(132, 11).Value = 15197.571
$ws.Cells.Item(132, 12).Value = 13232.7
$ws.Cells.Item(132, 13).Value = -12667.571
$ws.Cells.Item(132, 14).Value = -18292.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5240.76
$ws.Cells.Item(7, 9).Value = 3701.5
$ws.Cells.Item(7, 11).Value = 3701.5
$ws.Cells.Item(7, 13).Value = -3589.5
$ws.Cells.Item(46, 8).Value = 5454.5
$ws.Cells.Item(46, 10).Value = 4087.2632
$ws.Cells.Item(46, 12).Value = 4087.2632
$ws.Cells.Item(46, 14).Value = -4463.263199999999
$ws.Cells.Item(61, 8).Value = 1933.2084
$ws.Cells.Item(61, 9).Value = 1986.1389
$ws.Cells.Item(61, 10).Value = 1774.4166
$ws.Cells.Item(61, 11).Value = 1986.1389
$ws.Cells.Item(61, 12).Value = 1774.4166
$ws.Cells.Item(61, 13).Value = -1784.1389
$ws.Cells.Item(61, 14).Value = -2178.4166
$ws.Cells.Item(106, 8).Value = 2538886.8
$ws.Cells.Item(106, 10).Value = 2538886.8
$ws.Cells.Item(106, 12).Value = 2538886.8
$ws.Cells.Item(106, 14).Value = -2541410.8
$ws.Cells.Item(113, 8).Value = 1933.2084
$ws.Cells.Item(113, 9).Value = 1986.1389
$ws.Cells.Item(113, 10).Value = 1774.4166
$ws.Cells.Item(113, 11).Value = 1986.1389
$ws.Cells.Item(113, 12).Value = 1774.4166
$ws.Cells.Item(113, 13).Value = 183.8611000000001
$ws.Cells.Item(113, 14).Value = -6114.4166
$ws.Cells.Item(126, 8).Value = 5240.76
$ws.Cells.Item(126, 9).Value = 3701.5
$ws.Cells.Item(126, 11).Value = 11104.5
$ws.Cells.Item(126, 13).Value = -8634.5
$ws.Cells.Item(132, 8).Value = 304055.38
$ws.Cells.Item(132, 9).Value = 334212.5
$ws.Cells.Item(132, 11).Value = 1002637.5
$ws.Cells.Item(132, 13).Value = -1000107.5
$ws.Cells.Item(140, 8).Value = 419999
$ws.Cells.Item(140, 10).Value = 419999
$ws.Cells.Item(140, 12).Value = 419999
$ws.Cells.Item(140, 14).Value = -430359

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 24344
$ws.Cells.Item(20, 10).Value = 30011
$ws.Cells.Item(20, 12).Value = 30011
$ws.Cells.Item(20, 14).Value = -30491
$ws.Cells.Item(46, 8).Value = 69928
$ws.Cells.Item(46, 10).Value = 69928
$ws.Cells.Item(46, 12).Value = 69928
$ws.Cells.Item(46, 14).Value = -70390
$ws.Cells.Item(98, 8).Value = 42899
$ws.Cells.Item(98, 10).Value = 42899
$ws.Cells.Item(98, 12).Value = 42899
$ws.Cells.Item(98, 14).Value = -48889
$ws.Cells.Item(107, 8).Value = 555.2222
$ws.Cells.Item(107, 9).Value = 539.8333
$ws.Cells.Item(107, 11).Value = 1619.4999
$ws.Cells.Item(107, 13).Value = 300.5001
$ws.Cells.Item(122, 8).Value = 41671064
$ws.Cells.Item(122, 9).Value = 47622388
$ws.Cells.Item(122, 11).Value = 142867164
$ws.Cells.Item(122, 13).Value = -142864714
$ws.Cells.Item(132, 8).Value = 18448.184
$ws.Cells.Item(132, 9).Value = 22230.209
$ws.Cells.Item(132, 11).Value = 66690.62699999999
$ws.Cells.Item(132, 13).Value = -64160.62699999999
$ws.Cells.Item(134, 8).Value = 69928
$ws.Cells.Item(134, 10).Value = 69928
$ws.Cells.Item(134, 12).Value = 209784
$ws.Cells.Item(134, 14).Value = -214854
